# Apply cryptocurrency price/volume updates per commit:
# "Updated cryptos list on Mon Dec  4 08:38:21 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings (e.g. "233.20", "1.00")
# must be pre-formatted as Text so Excel keeps them as literal strings
# instead of silently converting them to floating-point numbers
# (which would corrupt values like trailing zeros / dotted thousands).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = '41.640.21'
$ws.Range("E2").Value = '  +5.48%  '
$ws.Range("D3").Value = '2.253.19'
$ws.Range("E3").Value = '  +4.11%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '233.20'
$ws.Range("E5").Value = '  +2.34%  '
$ws.Range("E6").Value = '  +3.29%  '
$ws.Range("D7").Value = '64.30'
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.412'
$ws.Range("E9").Value = '  +3.87%  '
$ws.Range("D10").Value = '60.18'
$ws.Range("E10").Value = '  +3.58%  '
$ws.Range("D11").Value = '0.0906'
$ws.Range("E11").Value = '  +5.99%  '
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("D13").Value = '2.587.27'
$ws.Range("E13").Value = '  +4.08%  '
$ws.Range("D14").Value = '16.19'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").Value = '22.72'
$ws.Range("E15").Value = '  +2.67%  '
$ws.Range("D16").Value = '0.831'
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("D17").Value = '5.69'
$ws.Range("E17").Value = '  +2.91%  '
$ws.Range("D18").Value = '2.250.48'
$ws.Range("E18").Value = '  +3.50%  '
$ws.Range("D19").Value = '41.508.83'
$ws.Range("E19").Value = '  +5.20%  '
$ws.Range("D20").Value = '0.0₃0938'
$ws.Range("E20").Value = '  +10.18%  '
$ws.Range("D21").Value = '73.87'
$ws.Range("E21").Value = '  +2.64%  '
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("D23").Value = '252.42'
$ws.Range("E23").Value = '  +9.96%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").Value = '9.86'
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("D28").Value = '0.149'
$ws.Range("E28").Value = '  +5.92%  '
$ws.Range("D29").Value = '173.17'
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("D30").Value = '20.49'
$ws.Range("E30").Value = '  +2.96%  '
$ws.Range("D31").Value = '1.45'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").Value = '2.82'
$ws.Range("E32").Value = '  +8.21%  '
$ws.Range("E33").Value = '  +2.84%  '
$ws.Range("D34").Value = '5.10'
$ws.Range("E34").Value = '  +7.96%  '
$ws.Range("E35").Value = '  +3.37%  '
$ws.Range("E36").Value = '  +3.77%  '
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").Value = '6.90'
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '3.87'
$ws.Range("E38").Value = '  +8.34%  '
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").Value = '0.000255'
$ws.Range("E40").Value = '  +61.27%  '
$ws.Range("D41").Value = '5.22'
$ws.Range("E41").Value = '  +22.53%  '
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '0.0242'
$ws.Range("E43").Value = '  +5.33%  '
$ws.Range("D44").Value = '8.87'
$ws.Range("E44").Value = '  +14.41%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '0.102'
$ws.Range("E45").Value = '  +9.48%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '102.88'
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '17.81'
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("E48").Value = '  +4.01%  '
$ws.Range("D49").Value = '1.512.07'
$ws.Range("E49").Value = '  -0.87%  '
$ws.Range("E50").Value = '  +3.14%  '
$ws.Range("E51").Value = '  -0.80%  '
